$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:G13")
$rng.NumberFormat = "0.00"
$rng.Select()
